$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 2566.3333
$ws.Range("I49").Value = 2000
$ws.Range("J49").Value = 2849.5
$ws.Range("K49").Value = 6000
$ws.Range("L49").Value = 8548.5
$ws.Range("M49").Value = -5864
$ws.Range("N49").Value = -8820.5

$ws.Range("H70").Value = 50152000
$ws.Range("I70").Value = 750000
$ws.Range("K70").Value = 2250000
$ws.Range("M70").Value = -2249730

$ws.Range("H73").Value = 50152000
$ws.Range("I73").Value = 750000
$ws.Range("K73").Value = 2250000
$ws.Range("M73").Value = -2249064

$ws.Range("H98").Value = 1983
$ws.Range("I98").Value = 1782.0857
$ws.Range("J98").Value = 5499
$ws.Range("K98").Value = 1782.0857
$ws.Range("L98").Value = 5499
$ws.Range("M98").Value = -284.0857000000001
$ws.Range("N98").Value = -8495

$ws.Range("H115").Value = 1006.8571
$ws.Range("I115").Value = 1194.25
$ws.Range("J115").Value = 757
$ws.Range("K115").Value = 3582.75
$ws.Range("L115").Value = 2271
$ws.Range("M115").Value = -2015.75
$ws.Range("N115").Value = -5405

$ws.Range("H116").Value = 41692780
$ws.Range("I116").Value = 55587150
$ws.Range("J116").Value = 9666
$ws.Range("K116").Value = 55587150
$ws.Range("L116").Value = 9666
$ws.Range("M116").Value = -55583708
$ws.Range("N116").Value = -16550

$ws.Range("H122").Value = 1983
$ws.Range("I122").Value = 1782.0857
$ws.Range("J122").Value = 5499
$ws.Range("K122").Value = 5346.257100000001
$ws.Range("L122").Value = 16497
$ws.Range("M122").Value = -2896.257100000001
$ws.Range("N122").Value = -21397

$ws.Range("H137").Value = 29759.807
$ws.Range("I137").Value = 41835.906
$ws.Range("K137").Value = 125507.718
$ws.Range("M137").Value = -122957.718

$ws.Range("H138").Value = 7356.7188
$ws.Range("J138").Value = 7895.4727
$ws.Range("L138").Value = 23686.4181
$ws.Range("N138").Value = -33966.4181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -2288

$ws.Range("H10").Value = 6683327.5
$ws.Range("J10").Value = 19993
$ws.Range("L10").Value = 19993
$ws.Range("N10").Value = -20333

$ws.Range("H32").Value = 1172354
$ws.Range("I32").Value = 1444120
$ws.Range("J32").Value = 39995.918
$ws.Range("K32").Value = 1444120
$ws.Range("L32").Value = 39995.918
$ws.Range("M32").Value = -1443833
$ws.Range("N32").Value = -40569.918

$ws.Range("H45").Value = 3693.1667
$ws.Range("I45").Value = 3991.8
$ws.Range("K45").Value = 3991.8
$ws.Range("M45").Value = -3614.8

$ws.Range("H102").Value = 5613.25
$ws.Range("I102").Value = 5613.25
$ws.Range("K102").Value = 5613.25
$ws.Range("M102").Value = -3991.25

$ws.Range("H122").Value = 6959.463
$ws.Range("I122").Value = 6668.294
$ws.Range("K122").Value = 20004.882
$ws.Range("M122").Value = -17554.882

$ws.Range("H132").Value = 315553.78
$ws.Range("I132").Value = 456714.8
$ws.Range("K132").Value = 1370144.4
$ws.Range("M132").Value = -1367614.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 110000
$ws.Range("J21").Value = 110000
$ws.Range("L21").Value = 110000
$ws.Range("N21").Value = -110472

$ws.Range("H107").Value = 3965.3333
$ws.Range("I107").Value = 3937.8
$ws.Range("K107").Value = 3937.8
$ws.Range("M107").Value = -2017.8

$ws.Range("H134").Value = 4820.2856
$ws.Range("I134").Value = 3652.6155
$ws.Range("K134").Value = 10957.8465
$ws.Range("M134").Value = -8422.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9135.451999999999
$ws.Range("I31").Value = 4657.773
$ws.Range("J31").Value = 14060.9
$ws.Range("K31").Value = 4657.773
$ws.Range("L31").Value = 14060.9
$ws.Range("M31").Value = -4362.773
$ws.Range("N31").Value = -14650.9

$ws.Range("H34").Value = 9135.451999999999
$ws.Range("I34").Value = 4657.773
$ws.Range("J34").Value = 14060.9
$ws.Range("K34").Value = 4657.773
$ws.Range("L34").Value = 14060.9
$ws.Range("M34").Value = -4455.773
$ws.Range("N34").Value = -14464.9

$ws.Range("H58").Value = 5214
$ws.Range("I58").Value = 2390
$ws.Range("K58").Value = 2390
$ws.Range("M58").Value = -2187

$ws.Range("H105").Value = 884.9
$ws.Range("I105").Value = 843.5
$ws.Range("K105").Value = 843.5
$ws.Range("M105").Value = 903.5

$ws.Range("H122").Value = 14292139
$ws.Range("I122").Value = 16667662
$ws.Range("K122").Value = 50002986
$ws.Range("M122").Value = -50000536

$ws.Range("H132").Value = 2916.8333
$ws.Range("I132").Value = 2606.4375
$ws.Range("K132").Value = 7819.3125
$ws.Range("M132").Value = -5289.3125

$ws.Range("H136").Value = 5214
$ws.Range("I136").Value = 2390
$ws.Range("K136").Value = 7170
$ws.Range("M136").Value = -4620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 4925
$ws.Range("I119").Value = 2566.6667
$ws.Range("K119").Value = 7700.000100000001
$ws.Range("M119").Value = -2862.000100000001

$ws.Range("H131").Value = 5999.75
$ws.Range("I131").Value = 10000
$ws.Range("K131").Value = 30000
$ws.Range("M131").Value = -24960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3833666.8

$ws.Range("H113").Value = 183
$ws.Range("I113").Value = 174.5
$ws.Range("K113").Value = 174.5
$ws.Range("M113").Value = 1995.5

$ws.Range("H122").Value = 2623.75
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400

$ws.Range("H132").Value = 8641.781000000001
$ws.Range("I132").Value = 8771.038
$ws.Range("K132").Value = 26313.114
$ws.Range("M132").Value = -23783.114

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6945.6665
$ws.Range("I7").Value = 4402.909
$ws.Range("K7").Value = 4402.909
$ws.Range("M7").Value = -4290.909

$ws.Range("H40").Value = 7759256
$ws.Range("I40").Value = 8553564
$ws.Range("J40").Value = 14750
$ws.Range("K40").Value = 8553564
$ws.Range("L40").Value = 14750
$ws.Range("M40").Value = -8553428
$ws.Range("N40").Value = -15022

$ws.Range("H61").Value = 5767.7036
$ws.Range("I61").Value = 6042.4585
$ws.Range("K61").Value = 6042.4585
$ws.Range("M61").Value = -5840.4585

$ws.Range("H113").Value = 5767.7036
$ws.Range("I113").Value = 6042.4585
$ws.Range("K113").Value = 6042.4585
$ws.Range("M113").Value = -3872.4585

$ws.Range("H126").Value = 6945.6665
$ws.Range("I126").Value = 4402.909
$ws.Range("K126").Value = 13208.727
$ws.Range("M126").Value = -10738.727

$ws.Range("H132").Value = 479827.75
$ws.Range("I132").Value = 717406.2
$ws.Range("J132").Value = 4670.857
$ws.Range("K132").Value = 2152218.6
$ws.Range("L132").Value = 14012.571
$ws.Range("M132").Value = -2149688.6
$ws.Range("N132").Value = -19072.571

$ws.Range("H136").Value = 8030.0835
$ws.Range("I136").Value = 7128.143
$ws.Range("K136").Value = 21384.429
$ws.Range("M136").Value = -18834.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 885.73334
$ws.Range("I96").Value = 777.63635
$ws.Range("J96").Value = 1183
$ws.Range("K96").Value = 777.63635
$ws.Range("L96").Value = 1183
$ws.Range("M96").Value = 595.36365
$ws.Range("N96").Value = -3929

$ws.Range("H122").Value = 31252216
$ws.Range("I122").Value = 52633230
$ws.Range("K122").Value = 157899690
$ws.Range("M122").Value = -157897240

$ws.Range("H126").Value = 8712.166999999999
$ws.Range("I126").Value = 8431.556
$ws.Range("K126").Value = 25294.668
$ws.Range("M126").Value = -22824.668

$ws.Range("H132").Value = 23547.527
$ws.Range("I132").Value = 33840.09
$ws.Range("K132").Value = 101520.27
$ws.Range("M132").Value = -98990.26999999999

$ws.Range("H136").Value = 76828.28999999999
$ws.Range("I136").Value = 4562.375
$ws.Range("J136").Value = 173182.83
$ws.Range("K136").Value = 13687.125
$ws.Range("L136").Value = 519548.49
$ws.Range("M136").Value = -11137.125
$ws.Range("N136").Value = -524648.49
